# Fruta / hortaliza, semanal
# Inserts two new weekly price observations (rows 283-284) into the
# "Brocoli" sheet, pushing the existing data for old rows 283..384 down
# to rows 285..386 (dimension grows from A1:R384 to A1:R386).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 283, shifting rows 283-384 down to 285-386.
$ws.Rows.Item(283).Resize(2).Insert()

# Copy static/common columns (A,B,C,E,F,G,H,N,Q,R) from the row that is now
# at 285 (the former row 283) into the two freshly inserted rows, then set
# the row-specific values (D,I,J,K,L,M,O,P) for each new row.

$commonCols = @(1,2,3,5,6,7,8,14,17,18)
foreach ($col in $commonCols) {
    $srcValue = $ws.Cells.Item(285, $col).Value2
    $ws.Cells.Item(283, $col).Value = $srcValue
    $ws.Cells.Item(284, $col).Value = $srcValue
}

# Ensure the date column keeps the same number format/style as the rest of
# column D (style index carrying the date/time format).
$ws.Cells.Item(283, 4).NumberFormat = $ws.Cells.Item(285, 4).NumberFormat
$ws.Cells.Item(284, 4).NumberFormat = $ws.Cells.Item(285, 4).NumberFormat

# New row 283 data
$ws.Cells.Item(283, 4).Value = 44468
$ws.Cells.Item(283, 9).Value = "Primera"
$ws.Cells.Item(283, 10).Value = 3400
$ws.Cells.Item(283, 11).Value = 600
$ws.Cells.Item(283, 12).Value = 650
$ws.Cells.Item(283, 13).Value = 625
$ws.Cells.Item(283, 15).Value = "Región Metropolitana"
$ws.Cells.Item(283, 16).Value = 625

# New row 284 data
$ws.Cells.Item(284, 4).Value = 44468
$ws.Cells.Item(284, 9).Value = "Segunda"
$ws.Cells.Item(284, 10).Value = 1600
$ws.Cells.Item(284, 11).Value = 450
$ws.Cells.Item(284, 12).Value = 500
$ws.Cells.Item(284, 13).Value = 475
$ws.Cells.Item(284, 15).Value = "Región Metropolitana"
$ws.Cells.Item(284, 16).Value = 475
